# Newly added iAuthor TC's - update generated credential values in row 2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Xlyzb626"
$ws.Range("B2").Value = 23111036
$ws.Range("C2").Value = "xgglbtq54"
$ws.Range("D2").Value = "AzB5%v`$3"
$ws.Range("F2").Value = "pUNmFVUM"
$ws.Range("G2").Value = "QXHF"
